$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep gridlines visible (avoid the engine's implicit showGridLines=false default on save)
$excel.ActiveWindow.DisplayGridlines = $true

# Header row values (becomes sharedStrings 0..3)
$ws.Range("A1").Value = "NUMERO DE CONTROL"
$ws.Range("B1").Value = "NOMBRE COMPLETO"
$ws.Range("C1").Value = "ESPECIALIDAD"
$ws.Range("D1").Value = "INSTITUCION"

# Bold + centered header row
$headerRange = $ws.Range("A1:D1")
$headerRange.HorizontalAlignment = -4108
$headerRange.Font.Bold = $true

# Column widths (A, B, C, D) - values chosen so the engine's pixel-rounded
# stored width lands as close as possible to the target (23.08 / 14.87 / 66.75)
$ws.Columns.Item(1).ColumnWidth = 22.16666667
$ws.Columns.Item(2).ColumnWidth = 22.16666667
$ws.Columns.Item(3).ColumnWidth = 14
$ws.Columns.Item(4).ColumnWidth = 65.83333333

# Touch G6 so it is included in the used range / dimension
$ws.Range("G6").NumberFormat = "General"

# Final selection lands on D6
$ws.Range("D6").Select()
